$d = $word.ActiveDocument

# ---- Phase 1: mark each moving text block with a unique placeholder (preserve trailing break to avoid run-merge) ----
$null = $d.Content.Find.Execute("Capacitar o aluno a identificar e gerenciar os riscos no ambiente de trabalho, com ênfase nos decorrentes das atividades em biotecnologia na indústria e em laboratórios de pesquisa.", $true, $false, $false, $false, $false, $true, 1, $false, "ZZPLACEHOLDER01", 2)
$null = $d.Content.Find.Execute("Empower the student to identify and manage risks in the workplace, with a focus on those arising from activities in biotechnology in industry and research laboratories.", $true, $false, $false, $false, $false, $true, 1, $false, "ZZPLACEHOLDER02", 2)
$null = $d.Content.Find.Execute("8711290 - Elisson Antônio da Costa Romanel^l", $true, $false, $false, $false, $false, $true, 1, $false, "ZZPLACEHOLDER03^l", 2)
$null = $d.Content.Find.Execute("8853480 - Tatiane da Franca Silva", $true, $false, $false, $false, $false, $true, 1, $false, "ZZPLACEHOLDER04", 2)
$null = $d.Content.Find.Execute("Introdução a segurança e medicina do trabalho; ^lConceitos técnico e aspectos legais em biossegurança;", $true, $false, $false, $false, $false, $true, 1, $false, "ZZPLACEHOLDER05", 2)
$null = $d.Content.Find.Execute("Introduction to occupational safety and medicine; Technical concepts and legal aspects in biosafety.", $true, $false, $false, $false, $false, $true, 1, $false, "ZZPLACEHOLDER06", 2)
$null = $d.Content.Find.Execute("Introdução a gestão em segurança do trabalho e estratégias de prevenção; Identificação e controle dos riscos ambientais (físicos, químicos e biológicos). Normas regulamentadoras;^lClasses de risco biológico, níveis de biossegurança e normas para a atividades de biotecnologia. ^lDescarte e classificação de resíduo;^lLegislação para produção e manejo organismos geneticamente modificados (OGM) e seus derivados;^lBiossegurança no manuseio de cobaias; ^lPrincípios de bioética;^lEstudos de casos problemas e soluções", $true, $false, $false, $false, $false, $true, 1, $false, "ZZPLACEHOLDER07", 2)
$null = $d.Content.Find.Execute("A avaliação será composta por provas, exercícios, projetos, seminários, relatórios e estudos de casos que poderão compor as notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.^l", $true, $false, $false, $false, $false, $true, 1, $false, "ZZPLACEHOLDER08^l", 2)
$null = $d.Content.Find.Execute("MF≥ 5,0 para aprovação 5,0. Prova de recuperação para alunos com 3,0≤MF<5,0^l", $true, $false, $false, $false, $false, $true, 1, $false, "ZZPLACEHOLDER09^l", 2)
$null = $d.Content.Find.Execute("(MF+RC)/2 ≥ 5,0 para aprovação, onde RC é uma prova de recuperação a ser aplicada", $true, $false, $false, $false, $false, $true, 1, $false, "ZZPLACEHOLDER10", 2)
$null = $d.Content.Find.Execute("1-Binsfeld, P. C. Fundamentos Técnicos e o Sistema Nacional de Biossegurança em Biotecnologia. Interciência, 1ª edição 2015.^l^l2-Gonçalves Simão, L. B. Gestão de Segurança e Medicina do Trabalho, Normas Regulamentadoras e Fator Acidentário de Prevenção. Cenofisco, 1ª edição 2015.^l^l3-Hirata, M.H., Mancini Filho, J. Hirata, R. D. C. Manual de biossegurança.  Editora Manole. 3ª edição 2016.^l^l4- Semplici, S. Onze Teses de Bioética. Editora Ideias e Letras;1ª edição 2014", $true, $false, $false, $false, $false, $true, 1, $false, "ZZPLACEHOLDER11", 2)

# ---- Phase 2: replace each placeholder (+ preserved break) with its final destination text ----
$null = $d.Content.Find.Execute("ZZPLACEHOLDER01", $true, $false, $false, $false, $false, $true, 1, $false, "Introdução a segurança e medicina do trabalho; ^lConceitos técnico e aspectos legais em biossegurança;", 2)
$null = $d.Content.Find.Execute("ZZPLACEHOLDER02", $true, $false, $false, $false, $false, $true, 1, $false, "Introduction to occupational safety and medicine; Technical concepts and legal aspects in biosafety.", 2)
$null = $d.Content.Find.Execute("ZZPLACEHOLDER03^l", $true, $false, $false, $false, $false, $true, 1, $false, "Capacitar o aluno a identificar e gerenciar os riscos no ambiente de trabalho, com ênfase nos decorrentes das atividades em biotecnologia na indústria e em laboratórios de pesquisa.^l", 2)
$null = $d.Content.Find.Execute("ZZPLACEHOLDER04", $true, $false, $false, $false, $false, $true, 1, $false, "Introdução a gestão em segurança do trabalho e estratégias de prevenção; Identificação e controle dos riscos ambientais (físicos, químicos e biológicos). Normas regulamentadoras;^lClasses de risco biológico, níveis de biossegurança e normas para a atividades de biotecnologia. ^lDescarte e classificação de resíduo;^lLegislação para produção e manejo organismos geneticamente modificados (OGM) e seus derivados;^lBiossegurança no manuseio de cobaias; ^lPrincípios de bioética;^lEstudos de casos problemas e soluções", 2)
$null = $d.Content.Find.Execute("ZZPLACEHOLDER05", $true, $false, $false, $false, $false, $true, 1, $false, "A avaliação será composta por provas, exercícios, projetos, seminários, relatórios e estudos de casos que poderão compor as notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.", 2)
$null = $d.Content.Find.Execute("ZZPLACEHOLDER06", $true, $false, $false, $false, $false, $true, 1, $false, "Empower the student to identify and manage risks in the workplace, with a focus on those arising from activities in biotechnology in industry and research laboratories.", 2)
$null = $d.Content.Find.Execute("ZZPLACEHOLDER07", $true, $false, $false, $false, $false, $true, 1, $false, "MF≥ 5,0 para aprovação 5,0. Prova de recuperação para alunos com 3,0≤MF<5,0", 2)
$null = $d.Content.Find.Execute("ZZPLACEHOLDER08^l", $true, $false, $false, $false, $false, $true, 1, $false, "(MF+RC)/2 ≥ 5,0 para aprovação, onde RC é uma prova de recuperação a ser aplicada^l", 2)
$null = $d.Content.Find.Execute("ZZPLACEHOLDER09^l", $true, $false, $false, $false, $false, $true, 1, $false, "1-Binsfeld, P. C. Fundamentos Técnicos e o Sistema Nacional de Biossegurança em Biotecnologia. Interciência, 1ª edição 2015.^l^l2-Gonçalves Simão, L. B. Gestão de Segurança e Medicina do Trabalho, Normas Regulamentadoras e Fator Acidentário de Prevenção. Cenofisco, 1ª edição 2015.^l^l3-Hirata, M.H., Mancini Filho, J. Hirata, R. D. C. Manual de biossegurança.  Editora Manole. 3ª edição 2016.^l^l4- Semplici, S. Onze Teses de Bioética. Editora Ideias e Letras;1ª edição 2014^l", 2)
$null = $d.Content.Find.Execute("ZZPLACEHOLDER10", $true, $false, $false, $false, $false, $true, 1, $false, "8711290 - Elisson Antônio da Costa Romanel", 2)
$null = $d.Content.Find.Execute("ZZPLACEHOLDER11", $true, $false, $false, $false, $false, $true, 1, $false, "8853480 - Tatiane da Franca Silva", 2)

Write-Host "done"